$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1177.654280618312
$ws.Range("B2").Value = 1182.785671819263
$ws.Range("C2").Value = 1160.771254458977
$ws.Range("D2").Value = 1262.987068965517
$ws.Range("E2").Value = 1776.650564803805

# Row 3
$ws.Range("A3").Value = 11.09120415449529
$ws.Range("B3").Value = 11.64729817708333
$ws.Range("C3").Value = 10.89588417114039
$ws.Range("D3").Value = 12.2727420402859
$ws.Range("E3").Value = 19.08367447191747

# Row 5
$ws.Range("A5").Value = 0.9206349206349206
$ws.Range("B5").Value = 0.9206349206349206
$ws.Range("C5").Value = 0.9206349206349206
$ws.Range("D5").Value = 0.9206349206349206
$ws.Range("E5").Value = 0.9206349206349206

# Row 6
$ws.Range("A6").Value = 0.001868703282417939
$ws.Range("B6").Value = 0.001860991989643175
$ws.Range("C6").Value = 0.001853194746595762
$ws.Range("D6").Value = 0.001859637774902975
$ws.Range("E6").Value = 0.00185648559205747
